$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07381966666666666
$ws.Range("H2").Value = 0.221459
$ws.Range("I2").Value = 0.1284640970637474
$ws.Range("J2").Value = 0.1284640970637474
$ws.Range("M2").Value = 4.902303000000001
$ws.Range("N2").Value = 14.706909
$ws.Range("O2").Value = 0.07597201518094217
$ws.Range("P2").Value = 0.07597201518094215
$ws.Range("Q2").Value = 0.361886373359
$ws.Range("R2").Value = 3.256977360231
$ws.Range("S2").Value = 0.009759676332333043
$ws.Range("T2").Value = 0.009759676332333041
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07381966666666666
$ws.Range("H3").Value = 0.221459
$ws.Range("I3").Value = 0.1284640970637474
$ws.Range("J3").Value = 0.1284640970637474
$ws.Range("O3").Value = 0.08309923851776384
$ws.Range("P3").Value = 0.08309923851776382
$ws.Range("Q3").Value = 0.395836308731111
$ws.Range("R3").Value = 3.56252677858
$ws.Range("S3").Value = 0.01067526864286951
$ws.Range("T3").Value = 0.01067526864286951
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07381966666666666
$ws.Range("H4").Value = 0.221459
$ws.Range("I4").Value = 0.1284640970637474
$ws.Range("J4").Value = 0.1284640970637474
$ws.Range("M4").Value = 30.73728233333334
$ws.Range("N4").Value = 92.21184700000001
$ws.Range("O4").Value = 0.476342094735659
$ws.Range("P4").Value = 0.4763420947356589
$ws.Range("Q4").Value = 2.269015936085889
$ws.Range("R4").Value = 20.421143424773
$ws.Range("S4").Value = 0.06119285709367044
$ws.Range("T4").Value = 0.06119285709367043
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07381966666666666
$ws.Range("H5").Value = 0.221459
$ws.Range("I5").Value = 0.1284640970637474
$ws.Range("J5").Value = 0.1284640970637474
$ws.Range("M5").Value = 3.761634666666666
$ws.Range("N5").Value = 11.284904
$ws.Range("O5").Value = 0.05829483938490915
$ws.Range("P5").Value = 0.05829483938490914
$ws.Range("Q5").Value = 0.2776826172151111
$ws.Range("R5").Value = 2.499143554935999
$ws.Range("S5").Value = 0.007488793905058531
$ws.Range("T5").Value = 0.00748879390505853
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.07381966666666666
$ws.Range("H6").Value = 0.221459
$ws.Range("I6").Value = 0.1284640970637474
$ws.Range("J6").Value = 0.1284640970637474
$ws.Range("M6").Value = 19.76432066666667
$ws.Range("N6").Value = 59.292962
$ws.Range("O6").Value = 0.3062918121807258
$ws.Range("P6").Value = 0.3062918121807258
$ws.Range("Q6").Value = 1.458995563506444
$ws.Range("R6").Value = 13.130960071558
$ws.Range("S6").Value = 0.03934750108981584
$ws.Range("T6").Value = 0.03934750108981584
$ws.Range("G7").Value = 0.1698756666666667
$ws.Range("H7").Value = 0.5096270000000001
$ws.Range("I7").Value = 0.2956247991470493
$ws.Range("J7").Value = 0.2956247991470493
$ws.Range("M7").Value = 4.902303000000001
$ws.Range("N7").Value = 14.706909
$ws.Range("O7").Value = 0.07597201518094217
$ws.Range("P7").Value = 0.07597201518094215
$ws.Range("Q7").Value = 0.8327819903270002
$ws.Range("R7").Value = 7.495037912943001
$ws.Range("S7").Value = 0.02245921172866261
$ws.Range("T7").Value = 0.02245921172866261
$ws.Range("G8").Value = 0.1698756666666667
$ws.Range("H8").Value = 0.5096270000000001
$ws.Range("I8").Value = 0.2956247991470493
$ws.Range("J8").Value = 0.2956247991470493
$ws.Range("O8").Value = 0.08309923851776384
$ws.Range("P8").Value = 0.08309923851776382
$ws.Range("Q8").Value = 0.9109084323044444
$ws.Range("R8").Value = 8.19817589074
$ws.Range("S8").Value = 0.02456619569608668
$ws.Range("T8").Value = 0.02456619569608667
$ws.Range("G9").Value = 0.1698756666666667
$ws.Range("H9").Value = 0.5096270000000001
$ws.Range("I9").Value = 0.2956247991470493
$ws.Range("J9").Value = 0.2956247991470493
$ws.Range("M9").Value = 30.73728233333334
$ws.Range("N9").Value = 92.21184700000001
$ws.Range("O9").Value = 0.476342094735659
$ws.Range("P9").Value = 0.4763420947356589
$ws.Range("Q9").Value = 5.221516327896556
$ws.Range("R9").Value = 46.99364695106901
$ws.Range("S9").Value = 0.1408185360815139
$ws.Range("T9").Value = 0.1408185360815139
$ws.Range("G10").Value = 0.1698756666666667
$ws.Range("H10").Value = 0.5096270000000001
$ws.Range("I10").Value = 0.2956247991470493
$ws.Range("J10").Value = 0.2956247991470493
$ws.Range("M10").Value = 3.761634666666666
$ws.Range("N10").Value = 11.284904
$ws.Range("O10").Value = 0.05829483938490915
$ws.Range("P10").Value = 0.05829483938490914
$ws.Range("Q10").Value = 0.6390101967564444
$ws.Range("R10").Value = 5.751091770808
$ws.Range("S10").Value = 0.01723340018447327
$ws.Range("T10").Value = 0.01723340018447327
$ws.Range("G11").Value = 0.1698756666666667
$ws.Range("H11").Value = 0.5096270000000001
$ws.Range("I11").Value = 0.2956247991470493
$ws.Range("J11").Value = 0.2956247991470493
$ws.Range("M11").Value = 19.76432066666667
$ws.Range("N11").Value = 59.292962
$ws.Range("O11").Value = 0.3062918121807258
$ws.Range("P11").Value = 0.3062918121807258
$ws.Range("Q11").Value = 3.357477149463778
$ws.Range("R11").Value = 30.217294345174
$ws.Range("S11").Value = 0.09054745545631282
$ws.Range("T11").Value = 0.09054745545631282
$ws.Range("G12").Value = 0.07389766666666667
$ws.Range("H12").Value = 0.221693
$ws.Range("I12").Value = 0.1285998359531712
$ws.Range("J12").Value = 0.1285998359531712
$ws.Range("M12").Value = 4.902303000000001
$ws.Range("N12").Value = 14.706909
$ws.Range("O12").Value = 0.07597201518094217
$ws.Range("P12").Value = 0.07597201518094215
$ws.Range("Q12").Value = 0.362268752993
$ws.Range("R12").Value = 3.260418776937
$ws.Range("S12").Value = 0.009769988689300997
$ws.Range("T12").Value = 0.009769988689300995
$ws.Range("G13").Value = 0.07389766666666667
$ws.Range("H13").Value = 0.221693
$ws.Range("I13").Value = 0.1285998359531712
$ws.Range("J13").Value = 0.1285998359531712
$ws.Range("O13").Value = 0.08309923851776384
$ws.Range("P13").Value = 0.08309923851776382
$ws.Range("Q13").Value = 0.3962545608511111
$ws.Range("R13").Value = 3.56629104766
$ws.Range("S13").Value = 0.01068654844121788
$ws.Range("T13").Value = 0.01068654844121788
$ws.Range("G14").Value = 0.07389766666666667
$ws.Range("H14").Value = 0.221693
$ws.Range("I14").Value = 0.1285998359531712
$ws.Range("J14").Value = 0.1285998359531712
$ws.Range("M14").Value = 30.73728233333334
$ws.Range("N14").Value = 92.21184700000001
$ws.Range("O14").Value = 0.476342094735659
$ws.Range("P14").Value = 0.4763420947356589
$ws.Range("Q14").Value = 2.271413444107889
$ws.Range("R14").Value = 20.442720996971
$ws.Range("S14").Value = 0.0612575152405957
$ws.Range("T14").Value = 0.06125751524059569
$ws.Range("G15").Value = 0.07389766666666667
$ws.Range("H15").Value = 0.221693
$ws.Range("I15").Value = 0.1285998359531712
$ws.Range("J15").Value = 0.1285998359531712
$ws.Range("M15").Value = 3.761634666666666
$ws.Range("N15").Value = 11.284904
$ws.Range("O15").Value = 0.05829483938490915
$ws.Range("P15").Value = 0.05829483938490914
$ws.Range("Q15").Value = 0.2779760247191111
$ws.Range("R15").Value = 2.501784222472
$ws.Range("S15").Value = 0.007496706781815782
$ws.Range("T15").Value = 0.007496706781815781
$ws.Range("G16").Value = 0.07389766666666667
$ws.Range("H16").Value = 0.221693
$ws.Range("I16").Value = 0.1285998359531712
$ws.Range("J16").Value = 0.1285998359531712
$ws.Range("M16").Value = 19.76432066666667
$ws.Range("N16").Value = 59.292962
$ws.Range("O16").Value = 0.3062918121807258
$ws.Range("P16").Value = 0.3062918121807258
$ws.Range("Q16").Value = 1.460537180518444
$ws.Range("R16").Value = 13.144834624666
$ws.Range("S16").Value = 0.03938907680024088
$ws.Range("T16").Value = 0.03938907680024088
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2570396666666667
$ws.Range("H17").Value = 0.771119
$ws.Range("I17").Value = 0.447311267836032
$ws.Range("J17").Value = 0.447311267836032
$ws.Range("M17").Value = 4.902303000000001
$ws.Range("N17").Value = 14.706909
$ws.Range("O17").Value = 0.07597201518094217
$ws.Range("P17").Value = 0.07597201518094215
$ws.Range("Q17").Value = 1.260086329019
$ws.Range("R17").Value = 11.340776961171
$ws.Range("S17").Value = 0.03398313843064552
$ws.Range("T17").Value = 0.03398313843064551
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.2570396666666667
$ws.Range("H18").Value = 0.771119
$ws.Range("I18").Value = 0.447311267836032
$ws.Range("J18").Value = 0.447311267836032
$ws.Range("O18").Value = 0.08309923851776384
$ws.Range("P18").Value = 0.08309923851776382
$ws.Range("Q18").Value = 1.378299814197778
$ws.Range("R18").Value = 12.40469832778
$ws.Range("S18").Value = 0.03717122573758977
$ws.Range("T18").Value = 0.03717122573758976
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.2570396666666667
$ws.Range("H19").Value = 0.771119
$ws.Range("I19").Value = 0.447311267836032
$ws.Range("J19").Value = 0.447311267836032
$ws.Range("M19").Value = 30.73728233333334
$ws.Range("N19").Value = 92.21184700000001
$ws.Range("O19").Value = 0.476342094735659
$ws.Range("P19").Value = 0.4763420947356589
$ws.Range("Q19").Value = 7.900700805199223
$ws.Range("R19").Value = 71.10630724679301
$ws.Range("S19").Value = 0.2130731863198789
$ws.Range("T19").Value = 0.2130731863198789
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.2570396666666667
$ws.Range("H20").Value = 0.771119
$ws.Range("I20").Value = 0.447311267836032
$ws.Range("J20").Value = 0.447311267836032
$ws.Range("M20").Value = 3.761634666666666
$ws.Range("N20").Value = 11.284904
$ws.Range("O20").Value = 0.05829483938490915
$ws.Range("P20").Value = 0.05829483938490914
$ws.Range("Q20").Value = 0.9668893208417777
$ws.Range("R20").Value = 8.702003887576
$ws.Range("S20").Value = 0.02607593851356156
$ws.Range("T20").Value = 0.02607593851356156
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.2570396666666667
$ws.Range("H21").Value = 0.771119
$ws.Range("I21").Value = 0.447311267836032
$ws.Range("J21").Value = 0.447311267836032
$ws.Range("M21").Value = 19.76432066666667
$ws.Range("N21").Value = 59.292962
$ws.Range("O21").Value = 0.3062918121807258
$ws.Range("P21").Value = 0.3062918121807258
$ws.Range("Q21").Value = 5.080214396053111
$ws.Range("R21").Value = 45.721929564478
$ws.Range("S21").Value = 0.1370077788343563
$ws.Range("T21").Value = 0.1370077788343563
